$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preserving its original
# "Normal" style (no numFmt / style index change) and inline/plain text type,
# even when the text looks like a number (e.g. "211.65").
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.902.13"
Set-TextValue "E2" "  +0.04%  "
Set-TextValue "D3" "1.635.11"
Set-TextValue "E3" "  +0.17%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "211.65"
Set-TextValue "E5" "  -0.04%  "
Set-TextValue "E6" "  -0.55%  "
Set-TextValue "E7" "  +0.05%  "
Set-TextValue "D8" "23.33"
Set-TextValue "E8" "  +0.42%  "
Set-TextValue "D9" "0.256"
Set-TextValue "E9" "  -0.43%  "
Set-TextValue "E10" "  -0.26%  "
Set-TextValue "E11" "  +0.40%  "
Set-TextValue "D12" "1.867.24"
Set-TextValue "E12" "  +0.17%  "
Set-TextValue "D13" "1.642.01"
Set-TextValue "E13" "  +0.53%  "
Set-TextValue "E14" "  -1.02%  "
Set-TextValue "E15" "  -0.67%  "
Set-TextValue "D16" "65.29"
Set-TextValue "E16" "  +0.18%  "
Set-TextValue "D17" "27.920.45"
Set-TextValue "D18" "228.67"
Set-TextValue "E18" "  -0.56%  "
Set-TextValue "E19" "  +2.72%  "
Set-TextValue "D20" "0.0₃0720"
Set-TextValue "E20" "  -0.22%  "
Set-TextValue "D22" "4.37"
Set-TextValue "E22" "  +0.06%  "
Set-TextValue "D23" "10.13"
Set-TextValue "E23" "  -2.20%  "
Set-TextValue "E24" "  +0.67%  "
Set-TextValue "D25" "156.03"
Set-TextValue "E25" "  +1.60%  "
Set-TextValue "E26" "  -0.39%  "
Set-TextValue "E27" "  -0.11%  "
Set-TextValue "D28" "15.54"
Set-TextValue "E28" "  -0.40%  "
Set-TextValue "E29" "  -0.06%  "
Set-TextValue "E30" "  +0.16%  "
Set-TextValue "E31" "  -0.08%  "
Set-TextValue "E32" "  +0.90%  "
Set-TextValue "E33" "  +1.25%  "
Set-TextValue "D34" "1.399.84"
Set-TextValue "E34" "  +0.33%  "
Set-TextValue "D35" "1.60"
Set-TextValue "E35" "  +2.81%  "
Set-TextValue "E36" "  +1.48%  "
Set-TextValue "E37" "  -0.72%  "
Set-TextValue "E38" "  +0.28%  "
Set-TextValue "D39" "0.559"
Set-TextValue "E39" "  -0.25%  "
Set-TextValue "D40" "0.850"
Set-TextValue "E40" "  -2.21%  "
Set-TextValue "E41" "  +0.05%  "
Set-TextValue "E42" "  -1.06%  "
Set-TextValue "E43" "  +2.42%  "
Set-TextValue "D44" "66.06"
Set-TextValue "E44" "  -1.22%  "
Set-TextValue "E45" "  -1.10%  "
Set-TextValue "D46" "1.776.04"
Set-TextValue "E46" "  +0.08%  "
Set-TextValue "D47" "2.14"
Set-TextValue "E47" "  -2.82%  "
Set-TextValue "D48" "88.64"
Set-TextValue "E48" "  +1.11%  "
Set-TextValue "D49" "0.102"
Set-TextValue "E50" "  -0.41%  "
Set-TextValue "D51" "7.63"
Set-TextValue "E51" "  +2.08%  "
